$d = $word.ActiveDocument

# The document ends with a single bulleted paragraph:
#   "...Otherwise, keep it going back to the welcome screen"
# followed immediately by the _GoBack bookmark.
#
# We need to split that into three bulleted paragraphs:
#   1) the existing sentence (unchanged)
#   2) a new sentence: "Go back and make all the titles refer to the
#      particular city or attraction they are looking at for each screen"
#   3) an empty paragraph that now carries the _GoBack bookmark
#
# Using Find/Replace with "^p" in the replacement text inserts real
# paragraph breaks (and correctly carries the trailing bookmark down into
# the newly-created final paragraph), while preserving the ListParagraph /
# numPr formatting of the paragraph being split.

[void]$d.Content.Find.Execute(
    "welcome screen",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "welcome screen^pGo back and make all the titles refer to the particular city or attraction they are looking at for each screen^p",
    2
)
